$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the "Datatype Complex" block (old row 32),
# pushing it down to make room for a new "str2str" method block.
$ws.Rows("32:35").Insert()

# Copy the row-formatting pattern from the existing "oneArg" method block
# (rows 28-31: blank-top-border / signature / body / blank-bottom-border)
# onto the newly inserted rows, since the new method block follows the
# exact same visual pattern. Copy one row at a time so Excel reuses the
# existing style/border definitions instead of synthesizing new ones.
$ws.Range("A28:E28").Copy()
$ws.Range("A32:E32").PasteSpecial(-4122)

$ws.Range("A29:E29").Copy()
$ws.Range("A33:E33").PasteSpecial(-4122)

$ws.Range("A30:E30").Copy()
$ws.Range("A34:E34").PasteSpecial(-4122)

$ws.Range("A31:E31").Copy()
$ws.Range("A35:E35").PasteSpecial(-4122)

# Merge the signature/body cells across B:C like the other method blocks.
# (Merging can nudge Excel into re-deriving border styles for the affected
# cells, so re-apply the original per-cell formatting afterwards.)
$ws.Range("B33:C33").Merge()
$ws.Range("B34:C34").Merge()

$ws.Range("A29:E29").Copy()
$ws.Range("A33:E33").PasteSpecial(-4122)

$ws.Range("A30:E30").Copy()
$ws.Range("A34:E34").PasteSpecial(-4122)

# Fill in the new method block's content.
$ws.Range("B33").Value = "Method String str2str(String data)"
$ws.Range("B34").Value = "return data;"
